# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they match the bold/bordered/centered look of
# the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-44) gets the same team record values.
$ws.Range("AD2:AD44").Value = 91
$ws.Range("AE2:AE44").Value = 71
$ws.Range("AF2:AF44").Value = 0
